$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shifted / swapped cell values in existing rows (per diff) ---
$ws.Range("D18").Value = 44462
$ws.Range("D19").Value = 44462
$ws.Range("D20").Value = 44222
$ws.Range("D21").Value = 44222
$ws.Range("D22").Value = 44159
$ws.Range("D23").Value = 44159
$ws.Range("D24").Value = 44330
$ws.Range("D25").Value = 44330
$ws.Range("D26").Value = 44231
$ws.Range("D27").Value = 44231
$ws.Range("D28").Value = 44348
$ws.Range("O28").Value = "Región de Ñuble"
$ws.Range("D29").Value = 44348
$ws.Range("I29").Value = "Segunda"
$ws.Range("J29").Value = 100
$ws.Range("K29").Value = 500
$ws.Range("L29").Value = 500
$ws.Range("M29").Value = 500
$ws.Range("P29").Value = 500
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 600
$ws.Range("L30").Value = 700
$ws.Range("M30").Value = 650
$ws.Range("P30").Value = 650
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 600
$ws.Range("L31").Value = 700
$ws.Range("M31").Value = 650
$ws.Range("P31").Value = 650
$ws.Range("D32").Value = 44435
$ws.Range("I32").Value = "Segunda"
$ws.Range("J32").Value = 100
$ws.Range("K32").Value = 500
$ws.Range("L32").Value = 500
$ws.Range("M32").Value = 500
$ws.Range("O32").Value = "Región Metropolitana"
$ws.Range("P32").Value = 500
$ws.Range("D33").Value = 44435
$ws.Range("D34").Value = 44224
$ws.Range("D35").Value = 44224
$ws.Range("D36").Value = 44313
$ws.Range("D37").Value = 44313
$ws.Range("D38").Value = 44169
$ws.Range("D39").Value = 44169
$ws.Range("D40").Value = 44271
$ws.Range("D41").Value = 44271
$ws.Range("D42").Value = 44322
$ws.Range("D43").Value = 44322
$ws.Range("D44").Value = 44355
$ws.Range("D45").Value = 44355
$ws.Range("D46").Value = 44299
$ws.Range("D47").Value = 44299
$ws.Range("D48").Value = 44320
$ws.Range("D49").Value = 44320
$ws.Range("D50").Value = 44391
$ws.Range("D51").Value = 44391
$ws.Range("D52").Value = 44229
$ws.Range("D53").Value = 44229
$ws.Range("D54").Value = 44245
$ws.Range("D55").Value = 44245
$ws.Range("D56").Value = 44334
$ws.Range("D57").Value = 44334
$ws.Range("D58").Value = 44250
$ws.Range("D59").Value = 44250
$ws.Range("D60").Value = 44434
$ws.Range("D61").Value = 44434
$ws.Range("D62").Value = 44327
$ws.Range("D63").Value = 44327
$ws.Range("D64").Value = 44383
$ws.Range("D65").Value = 44383
$ws.Range("D66").Value = 44405
$ws.Range("D67").Value = 44405
$ws.Range("D68").Value = 44161
$ws.Range("D69").Value = 44161
$ws.Range("D70").Value = 44336
$ws.Range("D71").Value = 44336
$ws.Range("D72").Value = 44341
$ws.Range("D73").Value = 44341
$ws.Range("D74").Value = 44274
$ws.Range("D75").Value = 44274
$ws.Range("D76").Value = 44400
$ws.Range("D77").Value = 44400
$ws.Range("D78").Value = 44442
$ws.Range("D79").Value = 44442
$ws.Range("D80").Value = 44453
$ws.Range("D81").Value = 44453
$ws.Range("D82").Value = 44208
$ws.Range("D83").Value = 44208
$ws.Range("D84").Value = 44420
$ws.Range("D85").Value = 44420
$ws.Range("D86").Value = 44237
$ws.Range("J86").Value = 200
$ws.Range("D87").Value = 44237
$ws.Range("J87").Value = 100
$ws.Range("D88").Value = 44285
$ws.Range("J88").Value = 100
$ws.Range("D89").Value = 44285
$ws.Range("J89").Value = 50
$ws.Range("D90").Value = 44344
$ws.Range("D91").Value = 44344
$ws.Range("D92").Value = 44217
$ws.Range("D93").Value = 44217
$ws.Range("D94").Value = 44266
$ws.Range("D95").Value = 44266
$ws.Range("D96").Value = 44350
$ws.Range("D97").Value = 44350
$ws.Range("D98").Value = 44460
$ws.Range("D99").Value = 44460
$ws.Range("D100").Value = 44427
$ws.Range("D101").Value = 44427
$ws.Range("D102").Value = 44280
$ws.Range("D103").Value = 44280
$ws.Range("D104").Value = 44447
$ws.Range("D105").Value = 44447
$ws.Range("D106").Value = 44267
$ws.Range("D107").Value = 44267
$ws.Range("D108").Value = 44187
$ws.Range("D109").Value = 44187
$ws.Range("D110").Value = 44386
$ws.Range("D111").Value = 44386
$ws.Range("D112").Value = 44308
$ws.Range("D113").Value = 44308
$ws.Range("D114").Value = 44264
$ws.Range("D115").Value = 44264
$ws.Range("D116").Value = 44196
$ws.Range("D117").Value = 44196
$ws.Range("D118").Value = 44243
$ws.Range("D119").Value = 44243
$ws.Range("K119").Value = 500
$ws.Range("L119").Value = 500
$ws.Range("M119").Value = 500
$ws.Range("P119").Value = 500
$ws.Range("D120").Value = 44252
$ws.Range("D121").Value = 44252
$ws.Range("K121").Value = 700
$ws.Range("L121").Value = 700
$ws.Range("M121").Value = 700
$ws.Range("P121").Value = 700
$ws.Range("D122").Value = 44166
$ws.Range("D123").Value = 44166
$ws.Range("D124").Value = 44168
$ws.Range("O124").Value = "Región de Ñuble"
$ws.Range("D125").Value = 44168
$ws.Range("O125").Value = "Región de Ñuble"
$ws.Range("D126").Value = 44433
$ws.Range("O126").Value = "Región Metropolitana"
$ws.Range("D127").Value = 44433
$ws.Range("O127").Value = "Región Metropolitana"
$ws.Range("D128").Value = 44371
$ws.Range("O128").Value = "Región de Ñuble"
$ws.Range("D129").Value = 44371
$ws.Range("O129").Value = "Región de Ñuble"
$ws.Range("D130").Value = 44316
$ws.Range("O130").Value = "Región Metropolitana"
$ws.Range("D131").Value = 44316
$ws.Range("O131").Value = "Región Metropolitana"
$ws.Range("D132").Value = 44273
$ws.Range("D133").Value = 44273
$ws.Range("D134").Value = 44209
$ws.Range("D135").Value = 44209
$ws.Range("D136").Value = 44365
$ws.Range("D137").Value = 44365
$ws.Range("D138").Value = 44306
$ws.Range("D139").Value = 44306
$ws.Range("D140").Value = 44215
$ws.Range("D141").Value = 44215
$ws.Range("D142").Value = 44257
$ws.Range("O142").Value = "Región de Ñuble"
$ws.Range("D143").Value = 44257
$ws.Range("O143").Value = "Región de Ñuble"
$ws.Range("D144").Value = 44239
$ws.Range("O144").Value = "Región Metropolitana"
$ws.Range("D145").Value = 44239
$ws.Range("O145").Value = "Región Metropolitana"
$ws.Range("D146").Value = 44376
$ws.Range("D147").Value = 44376
$ws.Range("D148").Value = 44292
$ws.Range("D149").Value = 44292
$ws.Range("D150").Value = 44358
$ws.Range("D151").Value = 44358
$ws.Range("D152").Value = 44211
$ws.Range("D153").Value = 44211

$ws.Range("A154").Value = 11
$ws.Range("B154").Value = "Vega Monumental Concepción"
$ws.Range("C154").Value = "Bíobío"
$ws.Range("D154").Value = 44425
$ws.Range("E154").Value = 8
$ws.Range("F154").Value = 100112009
$ws.Range("G154").Value = "Acelga"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 200
$ws.Range("K154").Value = 600
$ws.Range("L154").Value = 700
$ws.Range("M154").Value = 650
$ws.Range("N154").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O154").Value = "Región de Ñuble"
$ws.Range("P154").Value = 650
$ws.Range("Q154").Value = 1
$ws.Range("R154").Value = "Hortaliza"
$ws.Range("D154").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A155").Value = 11
$ws.Range("B155").Value = "Vega Monumental Concepción"
$ws.Range("C155").Value = "Bíobío"
$ws.Range("D155").Value = 44425
$ws.Range("E155").Value = 8
$ws.Range("F155").Value = 100112009
$ws.Range("G155").Value = "Acelga"
$ws.Range("H155").Value = "Sin especificar"
$ws.Range("I155").Value = "Segunda"
$ws.Range("J155").Value = 100
$ws.Range("K155").Value = 500
$ws.Range("L155").Value = 500
$ws.Range("M155").Value = 500
$ws.Range("N155").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O155").Value = "Región de Ñuble"
$ws.Range("P155").Value = 500
$ws.Range("Q155").Value = 1
$ws.Range("R155").Value = "Hortaliza"
$ws.Range("D155").NumberFormat = "YYYY-MM-DD HH:MM:SS"
